$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Deskcount")

# Deskcount updates
$ws.Range("C12").Value = 79    # Dallas - Cypress Waters
$ws.Range("C45").Value = 32    # Melbourne
$ws.Range("C46").Value = 561   # Mohali

# Include in Occupancy Calculation: flip these rows from Yes to No
$ws.Range("F17").Value = "No"  # Greenwood Village
$ws.Range("F22").Value = "No"  # Memphis
$ws.Range("F23").Value = "No"  # Brentwood
$ws.Range("F25").Value = "No"  # Newport Beach
$ws.Range("F38").Value = "No"  # Spokane
$ws.Range("F39").Value = "No"  # Tampa
$ws.Range("F48").Value = "No"  # Santiago
$ws.Range("F49").Value = "No"  # Sao Paulo
$ws.Range("F50").Value = "No"  # Singapore

# Selection state as last left by the editor
$ws.Range("C42").Select()
